$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1320.5385
$ws.Range("I15").Value = 1320.5385
$ws.Range("K15").Value = 3961.6155
$ws.Range("M15").Value = -3792.6155

$ws.Range("H17").Value = 5175000
$ws.Range("J17").Value = 5175000
$ws.Range("L17").Value = 15525000
$ws.Range("N17").Value = -15525336

$ws.Range("H34").Value = 1216.3334
$ws.Range("I34").Value = 1216.3334
$ws.Range("K34").Value = 1216.3334
$ws.Range("M34").Value = -1013.3334

$ws.Range("H36").Value = 1216.3334
$ws.Range("I36").Value = 1216.3334
$ws.Range("K36").Value = 1216.3334
$ws.Range("M36").Value = -501.3334

$ws.Range("H138").Value = 4092.9312
$ws.Range("J138").Value = 4399.9585
$ws.Range("L138").Value = 13199.8755
$ws.Range("N138").Value = -23479.8755

$ws.Range("H141").Value = 3850.818
$ws.Range("I141").Value = 3850.818
$ws.Range("K141").Value = 11552.454
$ws.Range("M141").Value = -6372.454000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21792.31
$ws.Range("I32").Value = 13079.32
$ws.Range("K32").Value = 13079.32
$ws.Range("M32").Value = -12792.32

$ws.Range("H37").Value = 23435
$ws.Range("J37").Value = 23435
$ws.Range("L37").Value = 23435
$ws.Range("N37").Value = -23981

$ws.Range("H44").Value = 34996.5
$ws.Range("J44").Value = 34996.5
$ws.Range("L44").Value = 34996.5
$ws.Range("N44").Value = -35972.5

$ws.Range("H55").Value = 24998
$ws.Range("J55").Value = 24998
$ws.Range("L55").Value = 24998
$ws.Range("N55").Value = -25628

$ws.Range("H61").Value = 4042.111
$ws.Range("I61").Value = 2998.2
$ws.Range("J61").Value = 5347
$ws.Range("K61").Value = 2998.2
$ws.Range("L61").Value = 5347
$ws.Range("M61").Value = -2786.2
$ws.Range("N61").Value = -5771

$ws.Range("H74").Value = 2232.2632
$ws.Range("I74").Value = 2142.5293
$ws.Range("K74").Value = 2142.5293
$ws.Range("M74").Value = -1268.5293

$ws.Range("H77").Value = 2232.2632
$ws.Range("I77").Value = 2142.5293
$ws.Range("K77").Value = 10712.6465
$ws.Range("M77").Value = -6344.646500000001

$ws.Range("H80").Value = 39999.168
$ws.Range("J80").Value = 39999.168
$ws.Range("L80").Value = 39999.168
$ws.Range("N80").Value = -41995.168

$ws.Range("H83").Value = 39999.168
$ws.Range("J83").Value = 39999.168
$ws.Range("L83").Value = 119997.504
$ws.Range("N83").Value = -129981.504

$ws.Range("H112").Value = 8000
$ws.Range("J112").Value = 8000
$ws.Range("L112").Value = 8000
$ws.Range("N112").Value = -10954

$ws.Range("H134").Value = 49999
$ws.Range("J134").Value = 49999
$ws.Range("L134").Value = 49999
$ws.Range("N134").Value = -60139

$ws.Range("H135").Value = 35000
$ws.Range("J135").Value = 35000
$ws.Range("L135").Value = 35000
$ws.Range("N135").Value = -45140

$ws.Range("H136").Value = 4042.111
$ws.Range("I136").Value = 2998.2
$ws.Range("J136").Value = 5347
$ws.Range("K136").Value = 8994.599999999999
$ws.Range("L136").Value = 16041
$ws.Range("M136").Value = -6444.599999999999
$ws.Range("N136").Value = -21141

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8333.333000000001
$ws.Range("I20").Value = 10000
$ws.Range("K20").Value = 10000
$ws.Range("M20").Value = -9753

$ws.Range("H35").Value = 15000
$ws.Range("J35").Value = 15000
$ws.Range("L35").Value = 15000
$ws.Range("N35").Value = -15620

$ws.Range("H86").Value = 1000
$ws.Range("J86").Value = 1000
$ws.Range("L86").Value = 1000
$ws.Range("N86").Value = -3246

$ws.Range("H89").Value = 1000
$ws.Range("J89").Value = 1000
$ws.Range("L89").Value = 5000
$ws.Range("N89").Value = -16232

$ws.Range("H134").Value = 14775.956
$ws.Range("I134").Value = 14141.167
$ws.Range("J134").Value = 15000
$ws.Range("K134").Value = 42423.501
$ws.Range("L134").Value = 45000
$ws.Range("M134").Value = -39888.501
$ws.Range("N134").Value = -50070

$ws.Range("H135").Value = 49999
$ws.Range("J135").Value = 49999
$ws.Range("L135").Value = 49999
$ws.Range("N135").Value = -60139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 485
$ws.Range("I25").Value = 485
$ws.Range("K25").Value = 485
$ws.Range("M25").Value = -311

$ws.Range("H31").Value = 5564
$ws.Range("I31").Value = 4558.143
$ws.Range("K31").Value = 4558.143
$ws.Range("M31").Value = -4263.143

$ws.Range("H34").Value = 5564
$ws.Range("I34").Value = 4558.143
$ws.Range("K34").Value = 4558.143
$ws.Range("M34").Value = -4356.143

$ws.Range("H41").Value = 17455.555
$ws.Range("J41").Value = 19271.428
$ws.Range("L41").Value = 19271.428
$ws.Range("N41").Value = -20127.428

$ws.Range("H107").Value = 834.3333
$ws.Range("I107").Value = 745
$ws.Range("J107").Value = 1013
$ws.Range("K107").Value = 745
$ws.Range("L107").Value = 1013
$ws.Range("M107").Value = 1175
$ws.Range("N107").Value = -4853

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 2000
$ws.Range("I112").Value = 2000
$ws.Range("K112").Value = 6000
$ws.Range("M112").Value = -4892

$ws.Range("H122").Value = 1519.4445
$ws.Range("J122").Value = 1551
$ws.Range("L122").Value = 13959
$ws.Range("N122").Value = -18859

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 11333.333
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

$ws.Range("H70").Value = 6614
$ws.Range("J70").Value = 7002.6
$ws.Range("L70").Value = 7002.6
$ws.Range("N70").Value = -7542.6

$ws.Range("H73").Value = 6614
$ws.Range("J73").Value = 7002.6
$ws.Range("L73").Value = 7002.6
$ws.Range("N73").Value = -8874.6

$ws.Range("H97").Value = 949.06665
$ws.Range("I97").Value = 1061.4166
$ws.Range("J97").Value = 499.66666
$ws.Range("K97").Value = 1061.4166
$ws.Range("L97").Value = 499.66666
$ws.Range("M97").Value = -565.4166
$ws.Range("N97").Value = -1491.66666

$ws.Range("H132").Value = 3747.25
$ws.Range("I132").Value = 2997
$ws.Range("J132").Value = 4497.5
$ws.Range("K132").Value = 8991
$ws.Range("L132").Value = 13492.5
$ws.Range("M132").Value = -6461
$ws.Range("N132").Value = -18552.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4651.75
$ws.Range("I40").Value = 4806.4546
$ws.Range("J40").Value = 2950
$ws.Range("K40").Value = 4806.4546
$ws.Range("L40").Value = 2950
$ws.Range("M40").Value = -4670.4546
$ws.Range("N40").Value = -3222

$ws.Range("H61").Value = 3674.8
$ws.Range("I61").Value = 3674.8
$ws.Range("K61").Value = 3674.8
$ws.Range("M61").Value = -3472.8

$ws.Range("H113").Value = 3674.8
$ws.Range("I113").Value = 3674.8
$ws.Range("K113").Value = 3674.8
$ws.Range("M113").Value = -1504.8

$ws.Range("H132").Value = 12398.571
$ws.Range("I132").Value = 12148.125
$ws.Range("J132").Value = 13200
$ws.Range("K132").Value = 36444.375
$ws.Range("L132").Value = 39600
$ws.Range("M132").Value = -33914.375
$ws.Range("N132").Value = -44660

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()

$ws.Range("H122").Value = 1669.6666

$ws.Range("H132").Value = 2671.5
$ws.Range("I132").Value = 2671.5
$ws.Range("K132").Value = 8014.5
$ws.Range("M132").Value = -5484.5
